# Auto-generated edit script: applies scheduled-runner value updates
# across ALC, ARM, CRP, CUL, LTW, WVR sheets in the Leve Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 4274016
$ws.Range("I92").Value = 5747485.5
$ws.Range("J92").Value = 953.4
$ws.Range("K92").Value = 5747485.5
$ws.Range("L92").Value = 953.4
$ws.Range("M92").Value = -5746237.5
$ws.Range("N92").Value = -3449.4
$ws.Range("H98").Value = 775.75
$ws.Range("I98").Value = 775.75
$ws.Range("K98").Value = 775.75
$ws.Range("M98").Value = 722.25
$ws.Range("H122").Value = 775.75
$ws.Range("I122").Value = 775.75
$ws.Range("K122").Value = 2327.25
$ws.Range("M122").Value = 122.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10808.438
$ws.Range("I32").Value = 3834.9016
$ws.Range("J32").Value = 26000.785
$ws.Range("K32").Value = 3834.9016
$ws.Range("L32").Value = 26000.785
$ws.Range("M32").Value = -3547.9016
$ws.Range("N32").Value = -26574.785
$ws.Range("H74").Value = 2263.46
$ws.Range("I74").Value = 1943.8823
$ws.Range("J74").Value = 2942.5625
$ws.Range("K74").Value = 1943.8823
$ws.Range("L74").Value = 2942.5625
$ws.Range("M74").Value = -1069.8823
$ws.Range("N74").Value = -4690.5625
$ws.Range("H77").Value = 2263.46
$ws.Range("I77").Value = 1943.8823
$ws.Range("J77").Value = 2942.5625
$ws.Range("K77").Value = 9719.4115
$ws.Range("L77").Value = 14712.8125
$ws.Range("M77").Value = -5351.4115
$ws.Range("N77").Value = -23448.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3723988
$ws.Range("I31").Value = 5433516.5
$ws.Range("J31").Value = 3249.8823
$ws.Range("K31").Value = 5433516.5
$ws.Range("L31").Value = 3249.8823
$ws.Range("M31").Value = -5433221.5
$ws.Range("N31").Value = -3839.8823
$ws.Range("H34").Value = 3723988
$ws.Range("I34").Value = 5433516.5
$ws.Range("J34").Value = 3249.8823
$ws.Range("K34").Value = 5433516.5
$ws.Range("L34").Value = 3249.8823
$ws.Range("M34").Value = -5433314.5
$ws.Range("N34").Value = -3653.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53.25
$ws.Range("I2").Value = 27.5
$ws.Range("J2").Value = 130.5
$ws.Range("K2").Value = 165
$ws.Range("L2").Value = 783
$ws.Range("M2").Value = -52
$ws.Range("N2").Value = -1009
$ws.Range("H5").Value = 755.11536
$ws.Range("I5").Value = 682.5238000000001
$ws.Range("K5").Value = 2047.5714
$ws.Range("M5").Value = -1935.5714
$ws.Range("H10").Value = 258.42856
$ws.Range("I10").Value = 241.8
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 725.4000000000001
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = -586.4000000000001
$ws.Range("N10").Value = -1178
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 15
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 159
$ws.Range("N19").ClearContents()
$ws.Range("H23").Value = 60.17647
$ws.Range("I23").Value = 29.375
$ws.Range("J23").Value = 87.55556
$ws.Range("K23").Value = 88.125
$ws.Range("L23").Value = 262.66668
$ws.Range("M23").Value = 146.875
$ws.Range("N23").Value = -732.66668
$ws.Range("H24").Value = 740
$ws.Range("J24").Value = 740
$ws.Range("L24").Value = 2220
$ws.Range("N24").Value = -2680
$ws.Range("H26").Value = 738
$ws.Range("I26").Value = 218.85715
$ws.Range("J26").Value = 1257.1428
$ws.Range("K26").Value = 656.5714499999999
$ws.Range("L26").Value = 3771.4284
$ws.Range("M26").Value = -368.5714499999999
$ws.Range("N26").Value = -4347.428400000001
$ws.Range("H29").Value = 559.5
$ws.Range("I29").Value = 44
$ws.Range("J29").Value = 927.7143
$ws.Range("K29").Value = 132
$ws.Range("L29").Value = 2783.1429
$ws.Range("M29").Value = 145
$ws.Range("N29").Value = -3337.1429
$ws.Range("H34").Value = 512.9231
$ws.Range("I34").Value = 58.923077
$ws.Range("J34").Value = 966.9231
$ws.Range("K34").Value = 176.769231
$ws.Range("L34").Value = 2900.7693
$ws.Range("M34").Value = -92.76923099999999
$ws.Range("N34").Value = -3068.7693
$ws.Range("H35").Value = 1875
$ws.Range("I35").Value = 500
$ws.Range("J35").Value = 6000
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = -1212
$ws.Range("N35").Value = -18576
$ws.Range("H36").Value = 598.1429000000001
$ws.Range("I36").Value = 567.4
$ws.Range("J36").Value = 675
$ws.Range("K36").Value = 1702.2
$ws.Range("L36").Value = 2025
$ws.Range("M36").Value = -1533.2
$ws.Range("N36").Value = -2363
$ws.Range("H40").Value = 5923.5884
$ws.Range("I40").Value = 41.125
$ws.Range("J40").Value = 11152.444
$ws.Range("K40").Value = 164.5
$ws.Range("L40").Value = 44609.776
$ws.Range("M40").Value = -95.5
$ws.Range("N40").Value = -44747.776
$ws.Range("H62").Value = 4399.615
$ws.Range("I62").Value = 900
$ws.Range("J62").Value = 4691.25
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 14073.75
$ws.Range("M62").Value = -2014
$ws.Range("N62").Value = -15445.75
$ws.Range("H65").Value = 4399.615
$ws.Range("I65").Value = 900
$ws.Range("J65").Value = 4691.25
$ws.Range("K65").Value = 8100
$ws.Range("L65").Value = 42221.25
$ws.Range("M65").Value = -4668
$ws.Range("N65").Value = -49085.25
$ws.Range("H97").Value = 495.2381
$ws.Range("I97").Value = 412.9
$ws.Range("J97").Value = 570.0909
$ws.Range("K97").Value = 1238.7
$ws.Range("L97").Value = 1710.2727
$ws.Range("M97").Value = -742.6999999999998
$ws.Range("N97").Value = -2702.2727
$ws.Range("H112").Value = 7933.278
$ws.Range("I112").Value = 1999.5
$ws.Range("J112").Value = 8675
$ws.Range("K112").Value = 5998.5
$ws.Range("L112").Value = 26025
$ws.Range("M112").Value = -4890.5
$ws.Range("N112").Value = -28241
$ws.Range("H122").Value = 500.65216
$ws.Range("J122").Value = 802.1429000000001
$ws.Range("L122").Value = 7219.2861
$ws.Range("N122").Value = -12119.2861
$ws.Range("H124").Value = 1908.75
$ws.Range("I124").Value = 1052
$ws.Range("J124").Value = 3336.6667
$ws.Range("K124").Value = 3156
$ws.Range("L124").Value = 10010.0001
$ws.Range("M124").Value = 1754
$ws.Range("N124").Value = -19830.0001
$ws.Range("H125").Value = 7463.6206
$ws.Range("I125").Value = 3574.8333
$ws.Range("J125").Value = 8478.087
$ws.Range("K125").Value = 10724.4999
$ws.Range("L125").Value = 25434.261
$ws.Range("M125").Value = -5804.499899999999
$ws.Range("N125").Value = -35274.261
$ws.Range("H126").Value = 2195.5334
$ws.Range("I126").Value = 500.5
$ws.Range("J126").Value = 2456.3076
$ws.Range("K126").Value = 1501.5
$ws.Range("L126").Value = 7368.9228
$ws.Range("M126").Value = 3438.5
$ws.Range("N126").Value = -17248.9228
$ws.Range("H129").Value = 1339.4615
$ws.Range("J129").Value = 1437
$ws.Range("L129").Value = 4311
$ws.Range("N129").Value = -14311
$ws.Range("H135").Value = 755.11536
$ws.Range("I135").Value = 682.5238000000001
$ws.Range("K135").Value = 6142.7142
$ws.Range("M135").Value = -3607.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 207.4762
$ws.Range("I55").Value = 241.3
$ws.Range("J55").Value = 176.72728
$ws.Range("K55").Value = 241.3
$ws.Range("L55").Value = 176.72728
$ws.Range("M55").Value = -68.30000000000001
$ws.Range("N55").Value = -522.7272800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 23628.572
$ws.Range("J106").Value = 23628.572
$ws.Range("L106").Value = 23628.572
$ws.Range("N106").Value = -26152.572
